$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "Betty"
$ws.Cells.Item(2,2).Value = "'57237345"
$ws.Cells.Item(2,2).ClearFormats()
$ws.Cells.Item(2,3).Value = "Cat"
$ws.Cells.Item(2,4).Value = "Cat Room G"
$ws.Cells.Item(2,5).Value = "'11/8/2024"
$ws.Cells.Item(2,5).ClearFormats()

$ws.Cells.Item(3,1).Value = "Pickle"
$ws.Cells.Item(3,2).Value = "'57091836"
$ws.Cells.Item(3,2).ClearFormats()
$ws.Cells.Item(3,3).Value = "Cat"
$ws.Cells.Item(3,4).Value = "Feature Room 1"
$ws.Cells.Item(3,5).Value = "'10/22/2024"
$ws.Cells.Item(3,5).ClearFormats()

$ws.Cells.Item(4,1).Value = "Crystal"
$ws.Cells.Item(4,2).Value = "'57091845"
$ws.Cells.Item(4,2).ClearFormats()
$ws.Cells.Item(4,3).Value = "Cat"
$ws.Cells.Item(4,4).Value = "Feature Room 1"
$ws.Cells.Item(4,5).Value = "'10/22/2024"
$ws.Cells.Item(4,5).ClearFormats()

$ws.Cells.Item(5,1).Value = "Colby Jack Cheese"
$ws.Cells.Item(5,2).Value = "'57670293"
$ws.Cells.Item(5,2).ClearFormats()
$ws.Cells.Item(5,3).Value = "Dog"
$ws.Cells.Item(5,4).Value = "Dog A"
$ws.Cells.Item(5,5).Value = "'1/17/2025"
$ws.Cells.Item(5,5).ClearFormats()

$ws.Cells.Item(6,1).Value = "Sabre"
$ws.Cells.Item(6,2).Value = "'57305889"
$ws.Cells.Item(6,2).ClearFormats()
$ws.Cells.Item(6,3).Value = "Cat"
$ws.Cells.Item(6,4).Value = "Foster Home"
$ws.Cells.Item(6,5).Value = "'11/21/2024"
$ws.Cells.Item(6,5).ClearFormats()

$ws.Cells.Item(7,1).Value = "Fleetwood"
$ws.Cells.Item(7,2).Value = "'57657154"
$ws.Cells.Item(7,2).ClearFormats()
$ws.Cells.Item(7,3).Value = "Cat"
$ws.Cells.Item(7,4).Value = "Condo Rooms"
$ws.Cells.Item(7,5).Value = "'1/16/2025"
$ws.Cells.Item(7,5).ClearFormats()

$ws.Cells.Item(8,1).Value = "Mac"
$ws.Cells.Item(8,2).Value = "'57657157"
$ws.Cells.Item(8,2).ClearFormats()
$ws.Cells.Item(8,3).Value = "Cat"
$ws.Cells.Item(8,4).Value = "Condo Rooms"
$ws.Cells.Item(8,5).Value = "'1/16/2025"
$ws.Cells.Item(8,5).ClearFormats()

$ws.Cells.Item(9,1).Value = "Bon Jovi"
$ws.Cells.Item(9,2).Value = "'57657185"
$ws.Cells.Item(9,2).ClearFormats()
$ws.Cells.Item(9,3).Value = "Cat"
$ws.Cells.Item(9,4).Value = "Condo Rooms"
$ws.Cells.Item(9,5).Value = "'1/16/2025"
$ws.Cells.Item(9,5).ClearFormats()

$ws.Cells.Item(10,1).Value = "Tigress"
$ws.Cells.Item(10,2).Value = "'57657181"
$ws.Cells.Item(10,2).ClearFormats()
$ws.Cells.Item(10,3).Value = "Cat"
$ws.Cells.Item(10,4).Value = "Foster Home"
$ws.Cells.Item(10,5).Value = "'1/16/2025"
$ws.Cells.Item(10,5).ClearFormats()

$ws.Cells.Item(11,1).Value = "Breadstick"
$ws.Cells.Item(11,2).Value = "'58011693"
$ws.Cells.Item(11,2).ClearFormats()
$ws.Cells.Item(11,3).Value = "Bird"
$ws.Cells.Item(11,4).Value = "Small Animals & Exotics"
$ws.Cells.Item(11,5).Value = "'3/13/2025"
$ws.Cells.Item(11,5).ClearFormats()

$ws.Cells.Item(12,1).Value = "Jane"
$ws.Cells.Item(12,2).Value = "'56455382"
$ws.Cells.Item(12,2).ClearFormats()
$ws.Cells.Item(12,3).Value = "Cat"
$ws.Cells.Item(12,4).Value = "Cat Room H"
$ws.Cells.Item(12,5).Value = "'1/29/2025"
$ws.Cells.Item(12,5).ClearFormats()

$ws.Cells.Item(13,1).Value = "Maraschino"
$ws.Cells.Item(13,2).Value = "'57954349"
$ws.Cells.Item(13,2).ClearFormats()
$ws.Cells.Item(13,3).Value = "Dog"
$ws.Cells.Item(13,4).Value = "Dog A"
$ws.Cells.Item(13,5).Value = "'3/4/2025"
$ws.Cells.Item(13,5).ClearFormats()

$ws.Cells.Item(14,1).Value = "Zeke"
$ws.Cells.Item(14,2).Value = "'57925440"
$ws.Cells.Item(14,2).ClearFormats()
$ws.Cells.Item(14,3).Value = "Dog"
$ws.Cells.Item(14,4).Value = "Dog B"
$ws.Cells.Item(14,5).Value = "'2/28/2025"
$ws.Cells.Item(14,5).ClearFormats()

$ws.Cells.Item(15,1).Value = "Whiskey"
$ws.Cells.Item(15,2).Value = "'57935094"
$ws.Cells.Item(15,2).ClearFormats()
$ws.Cells.Item(15,3).Value = "Dog"
$ws.Cells.Item(15,4).Value = "Dog F"
$ws.Cells.Item(15,5).Value = "'2/28/2025"
$ws.Cells.Item(15,5).ClearFormats()

$ws.Cells.Item(16,1).Value = "Montague"
$ws.Cells.Item(16,2).Value = "'57765213"
$ws.Cells.Item(16,2).ClearFormats()
$ws.Cells.Item(16,3).Value = "Cat"
$ws.Cells.Item(16,4).Value = "Foster Home"
$ws.Cells.Item(16,5).Value = "'2/1/2025"
$ws.Cells.Item(16,5).ClearFormats()

$ws.Cells.Item(17,1).Value = "Romeo"
$ws.Cells.Item(17,2).Value = "'57765219"
$ws.Cells.Item(17,2).ClearFormats()
$ws.Cells.Item(17,3).Value = "Cat"
$ws.Cells.Item(17,4).Value = "Foster Home"
$ws.Cells.Item(17,5).Value = "'2/1/2025"
$ws.Cells.Item(17,5).ClearFormats()

$ws.Cells.Item(18,1).Value = "Mike"
$ws.Cells.Item(18,2).Value = "'57945726"
$ws.Cells.Item(18,2).ClearFormats()
$ws.Cells.Item(18,3).Value = "Dog"
$ws.Cells.Item(18,4).Value = "Foster Home"
$ws.Cells.Item(18,5).Value = "'3/3/2025"
$ws.Cells.Item(18,5).ClearFormats()

$ws.Cells.Item(19,1).Value = "Ella"
$ws.Cells.Item(19,2).Value = "'57945729"
$ws.Cells.Item(19,2).ClearFormats()
$ws.Cells.Item(19,3).Value = "Dog"
$ws.Cells.Item(19,4).Value = "Foster Home"
$ws.Cells.Item(19,5).Value = "'3/3/2025"
$ws.Cells.Item(19,5).ClearFormats()

$ws.Cells.Item(20,1).Value = "Casher"
$ws.Cells.Item(20,2).Value = "'58057188"
$ws.Cells.Item(20,2).ClearFormats()
$ws.Cells.Item(20,3).Value = "Dog"
$ws.Cells.Item(20,4).Value = "Dog A"
$ws.Cells.Item(20,5).Value = "'3/21/2025"
$ws.Cells.Item(20,5).ClearFormats()
